$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You are the head of marketing analytics tasked with collecting customer feedback data for a retail company with 500 stores and 2 million annual customers. Your budget allows for surveying 10,000 customers. You are asked to maximize representativeness while maintaining cost-efficiency when collecting this data.Which method should you choose given the requirements and constraints?",
        "ques_type": 2,
        "options": [
            "Send email surveys to customers across different store locations, demographics, and purchase histories.",
            "Place feedback forms in 10,000 randomly selected shopping bags.",
            "Conduct in-person surveys at 50 randomly selected stores.",
            "Send email surveys to 10,000 randomly selected customers."
        ],
        "score": "Send email surveys to customers across different store locations, demographics, and purchase histories."
    },
    {
        "title": "You're a project manager evaluating two suppliers for components. Supplier A offers materials at $8 per unit with a 20% defect rate. Supplier B's materials cost $10 per unit with a 10% defect rate. You have a robust quality control process that ensures all defective parts can be identified and not used. The project aims to procure 100 non-defective components while minimizing costs and maintaining quality.Which of the following approaches should you choose?",
        "ques_type": 2,
        "options": [
            "Procure all components from Supplier A.",
            "Procure all components from Supplier B.",
            "Procure 63 units from A and 55 units from B.",
            "Procure 50 units from A and 50 units from B."
        ],
        "score": "Procure all components from Supplier A."
    },
    {
        "title": "You are the head of Sales, and your team is analyzing sales data from the past year, which includes monthly sales figures, marketing campaign expenditures, and economic indicators. Assuming no impact of any external unknown factors, the CEO wants a highly accurate prediction of future sales data that includes predicting revenue values for the next few quarters.What statistical method would be most appropriate for this predictive analysis?",
        "ques_type": 2,
        "options": [
            "Choose logistic regression to predict binary sales outcomes.",
            "Apply clustering algorithms for customer segmentation.",
            "Perform a causal inference analysis to identify the direct impact of specific marketing campaigns on sales.",
            "Use linear regression to model sales trends."
        ],
        "score": "Use linear regression to model sales trends."
    },
    {
        "title": "You're the product manager for a clothing store with an inventory turnover rate of four times per year. Your competitors maintain a higher rate, and you would like to increase your turnover rate to compete with them. What should be your primary action based on this key performance indicator (KPI)?",
        "ques_type": 2,
        "options": [
            "Focus on creating the right product mix.",
            "Calibrate product pricing.",
            "Increase product durability.",
            "Optimize inventory levels."
        ],
        "score": "Optimize inventory levels."
    }
]
'@

$ws.Rows(1).Delete()
$ws.Range("A1").Value = $newText
$ws.Range("A1").Style = "Normal"
$ws.Rows(1).RowHeight = $ws.StandardHeight
